$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 'Chinmay'
$ws.Range("C4").Value = 'Deolekar'
$ws.Range("D4").Value = 'chinmay.deolekar1@gmail.com'
$ws.Range("E4").Value = '$2b$10$0nOYJPyljeT.fDMLQCsIZuFFri9peb8bTApe2YbCoTN8QYy9u4DDS'
